$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '47.325.53'
$ws.Range('E2').Value = '  +3.50%  '
$ws.Range('D3').Value = '2.502.31'
$ws.Range('E3').Value = '  +2.74%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '324.38'
$ws.Range('E5').Value = '  +1.25%  '
$ws.Range('D6').Value = '109.40'
$ws.Range('E6').Value = '  +4.79%  '
$ws.Range('E7').Value = '  +1.66%  '
$ws.Range('D9').Value = '0.544'
$ws.Range('E9').Value = '  +1.96%  '
$ws.Range('D10').Value = '38.99'
$ws.Range('E10').Value = '  +8.52%  '
$ws.Range('D11').Value = '0.0817'
$ws.Range('E11').Value = '  +1.63%  '
$ws.Range('E12').Value = '  +1.26%  '
$ws.Range('D13').Value = '18.58'
$ws.Range('E13').Value = '  +1.54%  '
$ws.Range('E14').Value = '  +3.01%  '
$ws.Range('D15').Value = '2.888.74'
$ws.Range('E15').Value = '  +2.66%  '
$ws.Range('D16').Value = '2.503.45'
$ws.Range('E16').Value = '  +2.96%  '
$ws.Range('D17').Value = '0.859'
$ws.Range('E17').Value = '  +2.35%  '
$ws.Range('D18').Value = '47.288.45'
$ws.Range('E18').Value = '  +3.71%  '
$ws.Range('D19').Value = '13.02'
$ws.Range('E19').Value = '  +5.05%  '
$ws.Range('D20').Value = '6.73'
$ws.Range('E20').Value = '  +5.01%  '
$ws.Range('D21').Value = '0.0₃0948'
$ws.Range('E21').Value = '  +2.01%  '
$ws.Range('D22').Value = '71.07'
$ws.Range('E22').Value = '  -0.62%  '
$ws.Range('E23').Value = '  +7.85%  '
$ws.Range('D24').Value = '250.48'
$ws.Range('D25').Value = '2.60'
$ws.Range('E25').Value = '  +3.80%  '
$ws.Range('D26').Value = '26.21'
$ws.Range('E26').Value = '  +1.71%  '
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('E28').Value = '  +4.87%  '
$ws.Range('E29').Value = '  +4.00%  '
$ws.Range('D30').Value = '35.85'
$ws.Range('E30').Value = '  +6.95%  '
$ws.Range('D31').Value = '0.136'
$ws.Range('E31').Value = '  +6.11%  '
$ws.Range('D32').Value = '49.75'
$ws.Range('E32').Value = '  +0.76%  '
$ws.Range('D33').Value = '20.02'
$ws.Range('E33').Value = '  -0.72%  '
$ws.Range('E34').Value = '  +3.90%  '
$ws.Range('E35').Value = '  +4.40%  '
$ws.Range('E36').Value = '  +0.25%  '
$ws.Range('D37').Value = '4.79'
$ws.Range('E37').Value = '  +5.68%  '
$ws.Range('D38').Value = '2.00'
$ws.Range('E38').Value = '  +5.51%  '
$ws.Range('E39').Value = '  +3.37%  '
$ws.Range('E40').Value = '  +1.66%  '
$ws.Range('D41').Value = '122.69'
$ws.Range('E41').Value = '  -3.99%  '
$ws.Range('E42').Value = '  -1.82%  '
$ws.Range('D43').Value = '21.37'
$ws.Range('E43').Value = '  +2.62%  '
$ws.Range('E44').Value = '  +2.57%  '
$ws.Range('D45').Value = '1.987.72'
$ws.Range('E45').Value = '  +1.50%  '
$ws.Range('E46').Value = '  +2.99%  '
$ws.Range('E47').Value = '  -1.12%  '
$ws.Range('E48').Value = '  -1.46%  '
$ws.Range('D49').Value = '9.07'
$ws.Range('E49').Value = '  -0.74%  '
$ws.Range('D50').Value = '5.37'
$ws.Range('E50').Value = '  +10.97%  '
$ws.Range('D51').Value = '78.66'
$ws.Range('E51').Value = '  +1.54%  '
